# TC08 CDS Library source query fixed
# Replace the (broken) Participants-tab Cypher query with the corrected
# version that uses OPTIONAL MATCH so participants without matching
# genomic_info/library_source are still properly filtered out, and sorts
# the collected sample ids.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = @'
MATCH (p:participant)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
OPTIONAL MATCH (p)<--(diag:diagnosis)
OPTIONAL MATCH (samp)<--(f:file)
OPTIONAL MATCH (f)<--(g:genomic_info)
WITH s, p, samp, f, g, diag
WHERE g.library_source in ['VIRAL RNA']
WITH p
OPTIONAL MATCH (p)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
WITH s, p, apoc.coll.sort(collect(distinct samp.sample_id)) as samp
RETURN 
coalesce(p.participant_id,'') as `Participant ID`,
coalesce(s.study_name, '') as `Study Name`,
coalesce(s.phs_accession,'') as `Accession`,
coalesce(p.gender,'') as `Gender`,
coalesce(apoc.text.join(samp, ','), '') as `Samples`
ORDER BY p.participant_id
LIMIT 100
'@

# The new query text is longer, so the row needs to grow to keep showing
# it wrapped in full.
$ws.Rows.Item(2).RowHeight = 299.25

# Leave the selection on the query cells of the Participants row.
$ws.Range("B4:B9").Select()
